{"js": "// Update the worksheet date and the twelve three-digit-by-one-digit\n// multiplication problems to the new values from the latest generator run.\nconst replacements = [\n  [\"2025-10-15 Wednesday\", \"2025-10-16 Thursday\"],\n  [\"670\u00d78=\", \"572\u00d72=\"],\n  [\"435\u00d74=\", \"386\u00d76=\"],\n  [\"651\u00d72=\", \"545\u00d75=\"],\n  [\"218\u00d77=\", \"517\u00d79=\"],\n  [\"835\u00d75=\", \"703\u00d78=\"],\n  [\"265\u00d75=\", \"575\u00d74=\"],\n  [\"222\u00d77=\", \"351\u00d78=\"],\n  [\"477\u00d78=\", \"994\u00d72=\"],\n  [\"371\u00d78=\", \"901\u00d74=\"],\n  [\"486\u00d75=\", \"579\u00d78=\"],\n  [\"688\u00d78=\", \"568\u00d73=\"],\n  [\"655\u00d76=\", \"604\u00d72=\"],\n  [\"192\u00d75=\", \"562\u00d79=\"],\n  [\"454\u00d76=\", \"341\u00d79=\"],\n  [\"542\u00d77=\", \"151\u00d77=\"],\n  [\"786\u00d76=\", \"276\u00d75=\"],\n  [\"244\u00d74=\", \"158\u00d79=\"],\n  [\"892\u00d78=\", \"190\u00d74=\"],\n  [\"410\u00d74=\", \"979\u00d74=\"],\n  [\"446\u00d72=\", \"984\u00d79=\"],\n  [\"277\u00d72=\", \"880\u00d72=\"],\n  [\"236\u00d78=\", \"683\u00d78=\"],\n  [\"307\u00d72=\", \"461\u00d77=\"],\n  [\"335\u00d73=\", \"876\u00d77=\"],\n  [\"292\u00d77=\", \"428\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the twelve three-digit-by-one-digit\n# multiplication problems to the new values from the latest generator run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-15 Wednesday\", \"2025-10-16 Thursday\"),\n    @(\"670\u00d78=\", \"572\u00d72=\"),\n    @(\"435\u00d74=\", \"386\u00d76=\"),\n    @(\"651\u00d72=\", \"545\u00d75=\"),\n    @(\"218\u00d77=\", \"517\u00d79=\"),\n    @(\"835\u00d75=\", \"703\u00d78=\"),\n    @(\"265\u00d75=\", \"575\u00d74=\"),\n    @(\"222\u00d77=\", \"351\u00d78=\"),\n    @(\"477\u00d78=\", \"994\u00d72=\"),\n    @(\"371\u00d78=\", \"901\u00d74=\"),\n    @(\"486\u00d75=\", \"579\u00d78=\"),\n    @(\"688\u00d78=\", \"568\u00d73=\"),\n    @(\"655\u00d76=\", \"604\u00d72=\"),\n    @(\"192\u00d75=\", \"562\u00d79=\"),\n    @(\"454\u00d76=\", \"341\u00d79=\"),\n    @(\"542\u00d77=\", \"151\u00d77=\"),\n    @(\"786\u00d76=\", \"276\u00d75=\"),\n    @(\"244\u00d74=\", \"158\u00d79=\"),\n    @(\"892\u00d78=\", \"190\u00d74=\"),\n    @(\"410\u00d74=\", \"979\u00d74=\"),\n    @(\"446\u00d72=\", \"984\u00d79=\"),\n    @(\"277\u00d72=\", \"880\u00d72=\"),\n    @(\"236\u00d78=\", \"683\u00d78=\"),\n    @(\"307\u00d72=\", \"461\u00d77=\"),\n    @(\"335\u00d73=\", \"876\u00d77=\"),\n    @(\"292\u00d77=\", \"428\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
